# Added comparison to the workflow
# Insert a new "comparison" variable row into the metaware_meta_clean.csv
# codebook sheet, just above the existing "ref.r" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metaware_meta_clean.csv")

# Shift row 16 ("ref.r") and everything below it down by one row.
$ws.Rows("16:16").Insert()

$ws.Range("A16").Value = "comparison"
$ws.Range("B16").Value = "Indicator of whether effect size is indicative of a difference or difference-in-difference"
$ws.Range("C16").Value = "diff = difference `ndiff_diff = difference in difference"
